{"js": "// The word was not being accepted (ended on a rejecting state and could\n// not keep reading what was left in the input): trim the superfluous\n// trailing \"_\" from the first line and add a new \"aaa_\" entry below it.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the paragraph that currently reads \"b#b__\" and drop the extra\n// trailing underscore so it reads \"b#b_\".\nconst target = paragraphs.items.find((p) => p.text === \"b#b__\") || paragraphs.items[0];\n\ntarget.insertText(\"b#b_\", Word.InsertLocation.replace);\nawait context.sync();\n\n// Append a brand-new paragraph right after it with the extra test entry.\ntarget.insertParagraph(\"aaa_\", Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "# The word was not being accepted (ended on a rejecting state and could\n# not keep reading what was left in the input): trim the superfluous\n# trailing \"_\" from the first line and add a new \"aaa_\" entry below it.\n\n$d = $word.ActiveDocument\n\n# Locate the paragraph that currently reads \"b#b__\" and drop the extra\n# trailing underscore so it reads \"b#b_\".\n$targetIndex = 1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs($i).Range.Text.TrimEnd(\"`r\") -eq \"b#b__\") {\n        $targetIndex = $i\n        break\n    }\n}\n$target = $d.Paragraphs($targetIndex)\n$target.Range.Text = \"b#b_\"\n\n# Append a brand-new paragraph right after it with the extra test entry.\n$target.Range.InsertParagraphAfter()\n$d.Paragraphs($targetIndex + 1).Range.Text = \"aaa_\"\n"}
